# This script updates the "Orégano" hortaliza weekly price sheet.
# The underlying source rows (Fecha, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Precio $/Kg) have been reshuffled/resynced onto
# different date rows as part of the weekly data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 44335, 16, 10000, 10000, 10000, 3333),
    @(3, 44419, 16, 10000, 10000, 10000, 3333),
    @(4, 44307, 160, 10000, 10000, 10000, 3333),
    @(5, 44412, 25, 10000, 10500, 10260, 3420),
    @(6, 44433, 16, 10000, 10500, 10250, 3417),
    @(7, 44314, 16, 10000, 10000, 10000, 3333),
    @(8, 44503, 16, 8000, 9000, 8500, 2833),
    @(9, 44377, 16, 10000, 10500, 10250, 3417),
    @(10, 44293, 16, 10000, 10000, 10000, 3333),
    @(11, 44435, 16, 10000, 10500, 10250, 3417),
    @(12, 44356, 16, 10000, 10000, 10000, 3333),
    @(13, 44342, 17, 10000, 10000, 10000, 3333),
    @(14, 44524, 16, 9000, 10000, 9500, 3167),
    @(15, 44363, 16, 10000, 10000, 10000, 3333),
    @(16, 44517, 16, 9000, 10000, 9500, 3167),
    @(17, 44300, 16, 10000, 10000, 10000, 3333),
    @(18, 44426, 16, 10000, 10500, 10250, 3417),
    @(19, 44384, 25, 10000, 10500, 10260, 3420),
    @(20, 44349, 12, 10000, 10000, 10000, 3333),
    @(21, 44181, 10, 10000, 12000, 11000, 3667),
    @(22, 44482, 16, 9000, 10000, 9500, 3167),
    @(23, 44266, 160, 10000, 10000, 10000, 3333),
    @(24, 44539, 16, 9000, 10000, 9500, 3167),
    @(25, 44370, 16, 10000, 10500, 10250, 3417),
    @(26, 44475, 16, 9000, 10000, 9500, 3167),
    @(27, 44468, 16, 10000, 11000, 10500, 3500),
    @(28, 44489, 16, 9000, 10000, 9500, 3167),
    @(29, 44447, 16, 10000, 10500, 10250, 3417),
    @(30, 44328, 16, 10000, 10000, 10000, 3333),
    @(31, 44279, 16, 10000, 10000, 10000, 3333),
    @(32, 44391, 16, 10000, 10000, 10000, 3333),
    @(33, 44510, 16, 9000, 10000, 9500, 3167),
    @(34, 44321, 25, 10000, 10000, 10000, 3333),
    @(35, 44405, 16, 10000, 10500, 10250, 3417),
    @(36, 44175, 70, 12000, 12000, 12000, 4000),
    @(37, 44195, 30, 10000, 10000, 10000, 3333),
    @(38, 44398, 16, 10000, 10500, 10250, 3417),
    @(39, 44461, 16, 9500, 10000, 9750, 3250),
    @(40, 44454, 16, 9500, 10000, 9750, 3250),
    @(41, 44540, 32, 8500, 9000, 8719, 2906),
    @(42, 44272, 70, 10000, 10000, 10000, 3333)
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    $ws.Cells.Item($r, 4).Value  = $rowData[1]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $rowData[2]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $rowData[3]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $rowData[4]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $rowData[5]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $rowData[6]   # P - Precio $/Kg
}
